$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("B15").Value = 78713

# Row 16 / 17 - swap A, Q, R between the two rows
$ws.Range("A16").Value = 112241843
$ws.Range("Q16").Value = 554754
$ws.Range("R16").Value = 7006933

$ws.Range("A17").Value = 112241841
$ws.Range("Q17").Value = 554673
$ws.Range("R17").Value = 7006971

# Row 18
$ws.Range("B18").Value = 90113

# Row 32
$ws.Range("B32").Value = 78713

# Row 33
$ws.Range("A33").Value = 112241840
$ws.Range("B33").Value = 89553
$ws.Range("Q33").Value = 555033
$ws.Range("R33").Value = 7006894

# Row 34
$ws.Range("A34").Value = 112241839
$ws.Range("B34").Value = 89553
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 1202
$ws.Range("F34").Value = "Ullticka"
$ws.Range("G34").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H34").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q34").Value = 554840
$ws.Range("R34").Value = 7006980

# Row 35 (becomes the Tretåig hackspett / ringhack äldre record)
# Row 36 currently carries the empty K-N "age/sex/activity/method" cells that
# need to move over to row 35, so copy them across before anything else touches row 36.
$ws.Range("K36:N36").Copy($ws.Range("K35:N35"))

$ws.Range("A35").Value = 112241844
$ws.Range("B35").Value = 56430
$ws.Range("E35").Value = 100109
$ws.Range("F35").Value = "Tretåig hackspett"
$ws.Range("G35").Value = "Picoides tridactylus"
$ws.Range("H35").Value = "(Linnaeus, 1758)"
$ws.Range("Q35").Value = 554782
$ws.Range("R35").Value = 7006984
$ws.Range("AC35").Value = "ringhack äldre"

# Row 36 (becomes the Doftticka record)
$ws.Range("A36").Value = 112241866
$ws.Range("B36").Value = 90113
$ws.Range("D36").Value = "VU"
$ws.Range("E36").Value = 760
$ws.Range("F36").Value = "Doftticka"
$ws.Range("G36").Value = "Haploporus odorus"
$ws.Range("H36").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q36").Value = 554828
$ws.Range("R36").Value = 7006965
$ws.Range("K36:N36").ClearContents()
$ws.Range("AC36").ClearContents()
